$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainText($addr, $val) {
    $ws.Range($addr).Value = $val
}

function Set-NumericLookingText($addr, $val) {
    $ws.Range($addr).Formula = "'" + $val
}

Set-PlainText "D2" "68.553.71"
Set-PlainText "E2" "  +1.73%  "

Set-PlainText "D3" "3.779.17"
Set-PlainText "E3" "  +0.57%  "

Set-PlainText "E4" "  -0.15%  "

Set-NumericLookingText "D5" "596.83"
Set-PlainText "E5" "  +0.18%  "

Set-NumericLookingText "D6" "168.65"
Set-PlainText "E6" "  -0.02%  "

Set-PlainText "D7" "3.775.14"
Set-PlainText "E7" "  +0.54%  "

Set-PlainText "E8" "  -0.04%  "

Set-NumericLookingText "D9" "0.524"
Set-PlainText "E9" "  -0.96%  "

Set-NumericLookingText "D10" "0.163"
Set-PlainText "E10" "  -0.59%  "

Set-NumericLookingText "D11" "6.52"
Set-PlainText "E11" "  +0.76%  "

Set-NumericLookingText "D12" "0.450"
Set-PlainText "E12" "  -1.27%  "

Set-NumericLookingText "D13" "0.0000265"
Set-PlainText "E13" "  -1.60%  "

Set-NumericLookingText "D14" "36.69"
Set-PlainText "E14" "  +0.00%  "

Set-PlainText "D15" "4.414.40"
Set-PlainText "E15" "  +0.55%  "

Set-PlainText "D16" "3.779.05"
Set-PlainText "E16" "  +0.33%  "

Set-PlainText "D17" "68.520.65"
Set-PlainText "E17" "  +1.69%  "

Set-NumericLookingText "D18" "18.30"
Set-PlainText "E18" "  -3.10%  "

Set-NumericLookingText "D19" "7.07"
Set-PlainText "E19" "  -2.14%  "

Set-PlainText "E20" "  -0.23%  "

Set-NumericLookingText "D21" "10.93"
Set-PlainText "E21" "  +4.11%  "

Set-NumericLookingText "D22" "468.26"
Set-PlainText "E22" "  +0.19%  "

Set-NumericLookingText "D23" "0.705"
Set-PlainText "E23" "  -2.36%  "

Set-NumericLookingText "D24" "85.08"
Set-PlainText "E24" "  +1.47%  "

Set-NumericLookingText "D25" "0.0000144"
Set-PlainText "E25" "  -3.36%  "

Set-NumericLookingText "D26" "2.25"
Set-PlainText "E26" "  +1.02%  "

Set-NumericLookingText "D27" "12.22"
Set-PlainText "E27" "  +0.86%  "

Set-NumericLookingText "D28" "10.21"
Set-PlainText "E28" "  -0.67%  "

Set-PlainText "E29" "  +0.18%  "

Set-PlainText "D30" "3.927.22"
Set-PlainText "E30" "  +0.38%  "

Set-PlainText "E31" "  -3.37%  "

Set-NumericLookingText "D32" "7.43"
Set-PlainText "E32" "  -2.30%  "

Set-PlainText "B33" "EthereumClassic"
Set-PlainText "C33" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-NumericLookingText "D33" "30.16"
Set-PlainText "E33" "  -0.51%  "

Set-PlainText "B34" "ImmutableX"
Set-PlainText "C34" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-NumericLookingText "D34" "2.22"
Set-PlainText "E34" "  -0.97%  "

Set-NumericLookingText "D35" "9.31"
Set-PlainText "E35" "  +1.84%  "

Set-NumericLookingText "D36" "1.00"

Set-PlainText "D37" "3.733.53"
Set-PlainText "E37" "  +0.25%  "

Set-PlainText "E38" "  -2.64%  "

Set-NumericLookingText "D39" "3.47"
Set-PlainText "E39" "  -9.33%  "

Set-PlainText "E40" "  +1.67%  "

Set-PlainText "E41" "  +0.72%  "

Set-NumericLookingText "D42" "5.85"
Set-PlainText "E42" "  -0.45%  "

Set-NumericLookingText "D43" "0.999"
Set-PlainText "E43" "  -0.12%  "

Set-PlainText "B44" "TheGraph"
Set-PlainText "C44" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-NumericLookingText "D44" "0.309"
Set-PlainText "E44" "  -1.42%  "

Set-PlainText "B45" "USDe"
Set-PlainText "C45" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-NumericLookingText "D45" "1.00"
Set-PlainText "E45" "  +0.06%  "

Set-NumericLookingText "D46" "1.98"
Set-PlainText "E46" "  +1.43%  "

Set-NumericLookingText "D47" "44.00"
Set-PlainText "E47" "  +13.57%  "

Set-NumericLookingText "D48" "8.60"
Set-PlainText "E48" "  -1.20%  "

Set-NumericLookingText "D49" "410.72"
Set-PlainText "E49" "  +2.48%  "

Set-NumericLookingText "D50" "45.73"
Set-PlainText "E50" "  -1.04%  "

Set-NumericLookingText "D51" "145.72"
Set-PlainText "E51" "  +2.85%  "
